# Add season record columns (Wins, Losses, Ties) to the STL_1997 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
# Copy the formatting from the last existing header cell (AC1, which uses
# the bold/bordered/centered header style) onto the three new header cells
# so they reuse the same cell style instead of creating new ones.
$headerFormat = $ws.Range("AC1")
$headerFormat.Copy()

$wins = $ws.Range("AD1")
$losses = $ws.Range("AE1")
$ties = $ws.Range("AF1")

$wins.PasteSpecial(-4122)
$losses.PasteSpecial(-4122)
$ties.PasteSpecial(-4122)

$wins.Value = "Wins"
$losses.Value = "Losses"
$ties.Value = "Ties"

# ---- Data rows (rows 2-53) ----
# Every player row gets the team's 1997 season record: 73 wins, 89 losses,
# 0 ties.
for ($row = 2; $row -le 53; $row++) {
    $ws.Cells.Item($row, 30).Value = 73  # column AD = Wins
    $ws.Cells.Item($row, 31).Value = 89  # column AE = Losses
    $ws.Cells.Item($row, 32).Value = 0   # column AF = Ties
}
